# Generate Report for Handback
# Re-run of the handback report generator: the two source files tracked in
# this workbook were re-processed, producing fresh generated/correspond
# file names and timestamps. This script rewrites the three report sheets
# (Overview, zh-cn, de-de) to reflect the new run.

$wb = $excel.ActiveWorkbook

# ---- new identifiers coming out of the latest handback run ----
$oldFile1 = "5d24a5f0-095a-4f6c-ba13-767e5b80d782.md"
$newFile1 = "4df01f36-eb7a-4837-8ea3-295ed1881c43.md"
$oldFile2 = "84cb22ab-5658-4cb6-b7d7-fbe251bee46d.md"
$newFile2 = "ffff3d16f885-6cc2-4bb9-9e3b-459eb75109ac.md"

$newOverviewDate = "2016-08-31 19:12:46"

$newZhCnXlf = "4df01f36-eb7a-4837-8ea3-295ed1881c43.096e73f99289f61a7af2c0f49ccedd829c099418.zh-cn.xlf"
$newZhCnHandoffDate = "2016-08-31 19:12:40"
$newZhCnHandbackDate = "2016-08-31 19:13:14"

$newDeDeXlf = "4df01f36-eb7a-4837-8ea3-295ed1881c43.096e73f99289f61a7af2c0f49ccedd829c099418.de-de.xlf"
$newDeDeHandoffDate = "2016-08-31 19:12:46"
$newDeDeHandbackDate = "2016-08-31 19:13:21"

# original (unchanged) hyperlink target addresses, keyed by sheet+cell
$addrOverviewB2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6d992497c3ca85c678b12b62ef2857972c7b2f57/e2e/5d24a5f0-095a-4f6c-ba13-767e5b80d782.md"
$addrOverviewB3 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6d992497c3ca85c678b12b62ef2857972c7b2f57/e2e/84cb22ab-5658-4cb6-b7d7-fbe251bee46d.md"

$addrZhCnA2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6d992497c3ca85c678b12b62ef2857972c7b2f57/e2e/5d24a5f0-095a-4f6c-ba13-767e5b80d782.md"
$addrZhCnI2 = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/52e66d487c4a3f7f122961b3fed88fce088c293f/e2e/5d24a5f0-095a-4f6c-ba13-767e5b80d782.md"
$addrZhCnA3 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6d992497c3ca85c678b12b62ef2857972c7b2f57/e2e/84cb22ab-5658-4cb6-b7d7-fbe251bee46d.md"
$addrZhCnI3 = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/52e66d487c4a3f7f122961b3fed88fce088c293f/e2e/84cb22ab-5658-4cb6-b7d7-fbe251bee46d.md"

$addrDeDeA2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6d992497c3ca85c678b12b62ef2857972c7b2f57/e2e/5d24a5f0-095a-4f6c-ba13-767e5b80d782.md"
$addrDeDeI2 = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/2ef79684e6a2d23efc7eca840651660899e0013f/e2e/5d24a5f0-095a-4f6c-ba13-767e5b80d782.md"
$addrDeDeA3 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6d992497c3ca85c678b12b62ef2857972c7b2f57/e2e/84cb22ab-5658-4cb6-b7d7-fbe251bee46d.md"
$addrDeDeI3 = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/2ef79684e6a2d23efc7eca840651660899e0013f/e2e/84cb22ab-5658-4cb6-b7d7-fbe251bee46d.md"

# =====================================================================
# Sheet "Overview"
# =====================================================================
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFile1
$wsOverview.Range("A3").Value = $newFile2

$wsOverview.Range("B2").Value = "e2e\$newFile1"
$wsOverview.Range("B3").Value = "e2e\$newFile2"

$wsOverview.Range("G2").Value = $newOverviewDate
$wsOverview.Range("G3").Value = $newOverviewDate

# hyperlinks: this engine's Hyperlinks.Delete() clears the whole sheet
# collection regardless of which range it's invoked on, so clear once and
# rebuild every hyperlink on the sheet from scratch.
$wsOverview.Range("A1").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $addrOverviewB2, "", "", "e2e\$newFile1")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $addrOverviewB3, "", "", "e2e\$newFile2")

# =====================================================================
# Sheet "zh-cn"
# =====================================================================
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = $newFile1
$wsZhCn.Range("I2").Value = $newFile1
$wsZhCn.Range("A3").Value = $newFile2
$wsZhCn.Range("I3").Value = $newFile2

$wsZhCn.Range("G2").Value = $newZhCnXlf
$wsZhCn.Range("J2").Value = $newZhCnXlf
$wsZhCn.Range("G3").Value = $newZhCnXlf
$wsZhCn.Range("J3").Value = $newZhCnXlf

$wsZhCn.Range("H2").Value = $newZhCnHandoffDate
$wsZhCn.Range("H3").Value = $newZhCnHandoffDate

$wsZhCn.Range("K2").Value = $newZhCnHandbackDate
$wsZhCn.Range("K3").Value = $newZhCnHandbackDate

$wsZhCn.Range("A1").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $addrZhCnA2, "", "", $newFile1)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $addrZhCnI2, "", "", $newFile1)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $addrZhCnA3, "", "", $newFile2)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $addrZhCnI3, "", "", $newFile2)

# =====================================================================
# Sheet "de-de"
# =====================================================================
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = $newFile1
$wsDeDe.Range("I2").Value = $newFile1
$wsDeDe.Range("A3").Value = $newFile2
$wsDeDe.Range("I3").Value = $newFile2

$wsDeDe.Range("G2").Value = $newDeDeXlf
$wsDeDe.Range("J2").Value = $newDeDeXlf
$wsDeDe.Range("G3").Value = $newDeDeXlf
$wsDeDe.Range("J3").Value = $newDeDeXlf

$wsDeDe.Range("H2").Value = $newDeDeHandoffDate
$wsDeDe.Range("H3").Value = $newDeDeHandoffDate

$wsDeDe.Range("K2").Value = $newDeDeHandbackDate
$wsDeDe.Range("K3").Value = $newDeDeHandbackDate

$wsDeDe.Range("A1").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $addrDeDeA2, "", "", $newFile1)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $addrDeDeI2, "", "", $newFile1)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $addrDeDeA3, "", "", $newFile2)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $addrDeDeI3, "", "", $newFile2)

Write-Host "Handback report regenerated."
